# Update the "Latest Handoff Datetime" (column D) for the row 5 entry
# (65d69f13-...) on both the zh-cn and de-de status sheets, reflecting a
# fresh handoff generated for handback.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-20 03:08:07"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-20 03:08:17"
